$d = $word.ActiveDocument

# Locate the exact text "Version 1." in the document body and use its start
# offset as the anchor for every sub-range we touch below.
$findRange = $d.Content
$findRange.Find.Execute("Version 1.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$base = $findRange.Start

# --- Step 1: split the run "Version" into "Versi" + "on" ----------------
# A direct InsertBefore/InsertAfter in the middle of a run just grows that
# run's text; Word (and this host) only create a genuine run boundary when
# something structural - like a bookmark - starts/ends inside the run. So
# we drop a temporary bookmark across "on" and immediately delete it again;
# the bookmark is gone but the run split it forced stays behind, with no
# left-over run formatting.
$splitPoint = $d.Range($base + 5, $base + 7)
$d.Bookmarks.Add("__tmpSplit1", $splitPoint) | Out-Null
$d.Bookmarks("__tmpSplit1").Delete()

# --- Step 2: "1" -> "2" ---------------------------------------------------
$digit = $d.Range($base + 8, $base + 9)
$digit.Text = "2"

# --- Step 3: move the trailing "." so it lands after the _GoBack bookmark -
# Delete the "." that currently sits in the " 1." (now " 2.") run, then
# retype it immediately before the end of the paragraph - which, since the
# bookmark start/end pair is already there, places the new run after it.
$dot = $d.Range($base + 9, $base + 10)
$dot.Delete()

$para = $d.Range($base, $base).Paragraphs(1)
$paraEnd = $para.Range.End
$tail = $d.Range($paraEnd - 1, $paraEnd - 1)
$tail.InsertAfter(".")
